$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Datos actualizados" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 21:22"

# 2. Row 4 - Estados Unidos
$ws.Range("B4").Value = 330386
$ws.Range("C4").Value = 19029
$ws.Range("D4").Value = 17005
$ws.Range("E4").Value = 303937
$ws.Range("G4").Value = 992
$ws.Range("H4").Value = 9444

# 3. Row 7 - Alemania
$ws.Range("B7").Value = 100009
$ws.Range("C7").Value = 3917
$ws.Range("E7").Value = 69734
$ws.Range("G7").Value = 131
$ws.Range("H7").Value = 1575

# 4. Row 13 - Suiza
$ws.Range("E13").Value = 13970
$ws.Range("G13").Value = 49
$ws.Range("H13").Value = 715

# 5/6. Rows 36 & 37 - Japon and Pakistan swap ranking (Pakistan overtakes Japon)
#      Row 36 becomes Pakistan, Row 37 becomes Japon
$ws.Range("A36").Value = "Pakistan"
$ws.Range("B36").Value = 3157
$ws.Range("C36").Value = 339
$ws.Range("D36").Value = 211
$ws.Range("E36").Value = 2899
$ws.Range("F36").Value = 18
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 47

$ws.Range("A37").Value = "Japon"
$ws.Range("B37").Value = 3139
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 514
$ws.Range("E37").Value = 2548
$ws.Range("F37").Value = 64
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 77

# 7. Row 53
$ws.Range("D53").Value = 88
$ws.Range("E53").Value = 1362

# 8. Row 86
$ws.Range("D86").Value = 16
$ws.Range("E86").Value = 417
$ws.Range("F86").Value = 14

# 9. Row 145
$ws.Range("F145").Value = 0
